$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, shifting existing rows 44-48 down to 45-49.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the latest price report.
$ws.Range("A44").Value = 1
$ws.Range("B44").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C44").Value = "Arica y Parinacota"
$ws.Range("D44").Value = 44984
$ws.Range("E44").Value = 15
$ws.Range("F44").Value = "Fruta"
$ws.Range("G44").Value = 100104
$ws.Range("H44").Value = "Frutos de pepita"
$ws.Range("I44").Value = 100104005
$ws.Range("J44").Value = "Pera"
$ws.Range("K44").Value = "Packham's Triumph"
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 300
$ws.Range("N44").Value = 25000
$ws.Range("O44").Value = 26000
$ws.Range("P44").Value = 25500
$ws.Range("Q44").Value = "$/caja 20 kilos granel"
$ws.Range("R44").Value = "Región de O'Higgins"
$ws.Range("S44").Value = 1275
$ws.Range("T44").Value = 20
